# Update the NATMI ligand-receptor (Ntf3-Ntrk3) table with newly computed
# TPM-based values, and extend the table with the missing MuSCs-sending
# rows (rows 8-10) to complete the full 3x3 sending/target cluster matrix.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A..T
#  A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
#  E Ligand-expressing cells, F Ligand detection rate,
#  G Ligand average expression value, H Ligand total expression value,
#  I Ligand derived specificity of average expression value,
#  J Ligand derived specificity of total expression value,
#  K Receptor-expressing cells, L Receptor detection rate,
#  M Receptor average expression value, N Receptor total expression value,
#  O Receptor derived specificity of average expression value,
#  P Receptor derived specificity of total expression value,
#  Q Edge average expression weight, R Edge total expression weight,
#  S Edge average expression derived specificity,
#  T Edge total expression derived specificity
$data = @(
    @("ECs",   "Ntf3", "Ntrk3", "ECs",   3, 1, 5.305917,            15.917751,           0.4336744870332215, 0.4336744870332215, 1, 0.3333333333333333, 0.004481, 0.013443, 0.01111165297720057, 0.01111165297720057, 0.023775814077,      0.213982326693,      0.004818840404978626, 0.004818840404978627),
    @("ECs",   "Ntf3", "Ntrk3", "FAPs",  3, 1, 5.305917,            15.917751,           0.4336744870332215, 0.4336744870332215, 3, 1,                  0.1511493333333333, 0.453448, 0.3748089577628241, 0.3748089577628241, 0.801985817272,      7.217872355447999,   0.1625450824932491,   0.1625450824932491),
    @("ECs",   "Ntf3", "Ntrk3", "MuSCs", 3, 1, 5.305917,            15.917751,           0.4336744870332215, 0.4336744870332215, 3, 1,                  0.24764,             0.74292,  0.6140793892599753, 0.6140793892599754, 1.31395728588,       11.82561557292,      0.2663105641349938,   0.2663105641349938),
    @("FAPs",  "Ntf3", "Ntrk3", "ECs",   3, 1, 2.528563666666667,   7.585691000000001,  0.2066699405724794, 0.2066699405724794, 1, 0.3333333333333333, 0.004481, 0.013443, 0.01111165297720057, 0.01111165297720057, 0.01133049379033333, 0.101974444113,      0.002296444660460056, 0.002296444660460056),
    @("FAPs",  "Ntf3", "Ntrk3", "FAPs",  3, 1, 2.528563666666667,   7.585691000000001,  0.2066699405724794, 0.2066699405724794, 3, 1,                  0.1511493333333333, 0.453448, 0.3748089577628241, 0.3748089577628241, 0.3821907125075556,  3.439716412568,      0.07746174502687581,  0.07746174502687581),
    @("FAPs",  "Ntf3", "Ntrk3", "MuSCs", 3, 1, 2.528563666666667,   7.585691000000001,  0.2066699405724794, 0.2066699405724794, 3, 1,                  0.24764,             0.74292,  0.6140793892599753, 0.6140793892599754, 0.6261735064133334,  5.635561557720001,   0.1269117508851436,   0.1269117508851436),
    @("MuSCs", "Ntf3", "Ntrk3", "ECs",   3, 1, 4.400310999999999,  13.200933,           0.359655572394299,  0.359655572394299,  1, 0.3333333333333333, 0.004481, 0.013443, 0.01111165297720057, 0.01111165297720057, 0.019717793591,      0.177460142319,      0.003996367911761888, 0.003996367911761889),
    @("MuSCs", "Ntf3", "Ntrk3", "FAPs",  3, 1, 4.400310999999999,  13.200933,           0.359655572394299,  0.359655572394299,  3, 1,                  0.1511493333333333, 0.453448, 0.3748089577628241, 0.3748089577628241, 0.6651040741093333,  5.985936666983999,   0.1348021302426991,   0.1348021302426991),
    @("MuSCs", "Ntf3", "Ntrk3", "MuSCs", 3, 1, 4.400310999999999,  13.200933,           0.359655572394299,  0.359655572394299,  3, 1,                  0.24764,             0.74292,  0.6140793892599753, 0.6140793892599754, 1.08969301604,       9.80723714436,       0.220857074239838,    0.220857074239838)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowVals[$c]
    }
}
